$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Type Of Items" column (column K) entirely. Everything to its
# right (Service Type .. Remarks) shifts one column left, and the now-unused
# "Type Of Items" shared string drops out of the table.
$ws.Columns.Item(11).Delete()

# Match the cursor position saved with the new layout.
$ws.Range("J8").Select()
